# Update the availability status of the "MacBook Air" row (row 3) on the
# productInfoData sheet from "In Stock" to "Out Of Stock".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("productInfoData")
$ws.Range("D3").Value = "Out Of Stock"
